# Apply the "fixed trendlines" edit to Sheet1.
#
# Summary of the change (see commit message: "I fixed most of the
# trendlines. The Log trendline still needs work."):
#   - E6 on Sheet1 now computes Sum(x^n) generically using an exponent
#     held in column F, instead of the old hard-coded Sum(x^2) formula.
#   - New helper rows 7-10 are added with labels Sumx^3..Sumx^6 and
#     exponents 3..6 in column F (the actual summation formulas for
#     those rows were not filled in yet - matches "still needs work").
#   - The old matrix-inversion scratch area (rows 9-16 in columns C-F)
#     is removed.
#   - The second "solve the system of equations" block (rows 20-30),
#     which is duplicated verbatim on Sheet2, is stripped down to just
#     its formatted-but-empty shell on Sheet1 (labels/values removed,
#     but the styled cells A20:B27 remain so the formatting survives).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 6: generalize the Sum(x^2) formula, add exponent helper cell ---
$ws.Range("E6").Formula = "=`$A`$2^F6+`$A`$3^F6+`$A`$4^F6+`$A`$5^F6"
$ws.Range("F6").Value = 2

# --- New rows 7-10: additional power-sum helper rows ---
$ws.Range("D7").Value = "Sumx^3"
$ws.Range("F7").Value = 3

$ws.Range("D8").Value = "Sumx^4"
$ws.Range("F8").Value = 4

# Row 9 previously held the MINVERSE scratch (C9:F9) - replace it.
$ws.Range("C9:F9").ClearContents()
$ws.Range("D9").Value = "Sumx^5"
$ws.Range("F9").Value = 5

# Row 10 previously held the MINVERSE scratch (C10:F10) - replace it.
$ws.Range("C10:F10").ClearContents()
$ws.Range("D10").Value = "Sumx^6"
$ws.Range("F10").Value = 6

# --- Remove the old matrix-inversion-for-a,b scratch area ---
$ws.Range("B12:F13").ClearContents()
$ws.Range("C15:C16").ClearContents()

# --- Strip the duplicated "system of equations" block on Sheet1 down to ---
# --- its empty, formatted shell (rows 20-30), leaving rows 20,21,22,23,24,27 ---
$ws.Range("A20:M30").ClearContents()

# Rows 22-24 lose their larger-font matrix labels, so their auto row
# height shrinks back down to the default text height (matches rows
# 20/27, which were already at that height).
$ws.Rows("22:24").RowHeight = 15.75

# --- Sheet view / selection bookkeeping ---
$ws.Range("E6").Select()

$wb.Save()
